# Daily attendance processing - 2026-01-19 23:35:58
# Re-sort the "Recorded By" (column G) list of names/emails in ascending
# ordinal (case-sensitive, uppercase-before-lowercase) order for every
# data row on the active sheet.

function Compare-Ordinal($a, $b) {
    $ordLenA = $a.Length
    $ordLenB = $b.Length
    $ordMinLen = [Math]::Min($ordLenA, $ordLenB)
    for ($ordIdx = 0; $ordIdx -lt $ordMinLen; $ordIdx++) {
        $ordCa = [int][char]$a[$ordIdx]
        $ordCb = [int][char]$b[$ordIdx]
        if ($ordCa -lt $ordCb) { return -1 }
        if ($ordCa -gt $ordCb) { return 1 }
    }
    if ($ordLenA -lt $ordLenB) { return -1 }
    if ($ordLenA -gt $ordLenB) { return 1 }
    return 0
}

function Sort-Ordinal($list) {
    $sortN = $list.Count
    for ($sortI = 0; $sortI -lt $sortN; $sortI++) {
        for ($sortJ = 0; $sortJ -lt ($sortN - $sortI - 1); $sortJ++) {
            $sortCmp = Compare-Ordinal $list[$sortJ] $list[$sortJ+1]
            if ($sortCmp -gt 0) {
                $sortTmp = $list[$sortJ]
                $list[$sortJ] = $list[$sortJ+1]
                $list[$sortJ+1] = $sortTmp
            }
        }
    }
    return $list
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($rowNum = 2; $rowNum -le $lastRow; $rowNum++) {
    $cell = $ws.Cells.Item($rowNum, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        $trimmedParts = @()
        foreach ($part in $parts) {
            $trimmedParts += $part.Trim()
        }

        $sortedParts = Sort-Ordinal $trimmedParts
        $newVal = [string]::Join(", ", $sortedParts)

        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
